$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of "Alcachofa" price records (fecha 45147) needs to be inserted
# right above the existing block that starts at row 1013. Inserting 3 whole
# rows there pushes every following row (old 1013..1059) down by 3, which
# matches the diff (old 1013-1015 reappear as new 1016-1018, ... old
# 1057-1059 reappear as new 1060-1062) and grows the used range to R1062.
$ws.Rows.Item(1013).Insert()
$ws.Rows.Item(1013).Insert()
$ws.Rows.Item(1013).Insert()

# Row 1013 - Argentina(o) / Primera
$ws.Range("A1013").Value = 6
$ws.Range("B1013").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1013").Value = "Metropolitana"
$ws.Range("D1013").Value = 45147
$ws.Range("E1013").Value = 13
$ws.Range("F1013").Value = 100112013
$ws.Range("G1013").Value = "Alcachofa"
$ws.Range("H1013").Value = "Argentina(o)"
$ws.Range("I1013").Value = "Primera"
$ws.Range("J1013").Value = 390
$ws.Range("K1013").Value = 11000
$ws.Range("L1013").Value = 12000
$ws.Range("M1013").Value = 11564
$ws.Range("N1013").Value = "$/caja 50 unidades"
$ws.Range("O1013").Value = "Provincia de Limarí"
$ws.Range("P1013").Value = 231
$ws.Range("Q1013").Value = 50
$ws.Range("R1013").Value = "Hortaliza"

# Row 1014 - Española / Extra
$ws.Range("A1014").Value = 6
$ws.Range("B1014").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1014").Value = "Metropolitana"
$ws.Range("D1014").Value = 45147
$ws.Range("E1014").Value = 13
$ws.Range("F1014").Value = 100112013
$ws.Range("G1014").Value = "Alcachofa"
$ws.Range("H1014").Value = "Española"
$ws.Range("I1014").Value = "Extra"
$ws.Range("J1014").Value = 540
$ws.Range("K1014").Value = 12000
$ws.Range("L1014").Value = 13000
$ws.Range("M1014").Value = 12500
$ws.Range("N1014").Value = "$/caja 25 unidades"
$ws.Range("O1014").Value = "Provincia de Limarí"
$ws.Range("P1014").Value = 12500
$ws.Range("Q1014").Value = 1
$ws.Range("R1014").Value = "Hortaliza"

# Row 1015 - Española / Primera
$ws.Range("A1015").Value = 6
$ws.Range("B1015").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1015").Value = "Metropolitana"
$ws.Range("D1015").Value = 45147
$ws.Range("E1015").Value = 13
$ws.Range("F1015").Value = 100112013
$ws.Range("G1015").Value = "Alcachofa"
$ws.Range("H1015").Value = "Española"
$ws.Range("I1015").Value = "Primera"
$ws.Range("J1015").Value = 550
$ws.Range("K1015").Value = 13000
$ws.Range("L1015").Value = 14000
$ws.Range("M1015").Value = 13545
$ws.Range("N1015").Value = "$/caja 30 unidades"
$ws.Range("O1015").Value = "Provincia de Limarí"
$ws.Range("P1015").Value = 452
$ws.Range("Q1015").Value = 30
$ws.Range("R1015").Value = "Hortaliza"
